$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.1046505537244214
$ws.Range("D2").Value = 0.9176013430670817

$ws.Range("C3").Value = -0.4316695768059721
$ws.Range("D3").Value = 0.6701825054794959

$ws.Range("C4").Value = -0.9368278429612501
$ws.Range("D4").Value = 0.3590182584847943

$ws.Range("C5").Value = -0.989650222052978
$ws.Range("D5").Value = 0.3331055904657108

$ws.Range("C6").Value = -0.1831442341051032
$ws.Range("D6").Value = 0.8563622053368978

$ws.Range("C7").Value = -0.5494661824306761
$ws.Range("D7").Value = 0.5882203648072881

$ws.Range("C8").Value = -1.108290756099964
$ws.Range("D8").Value = 0.2797060300312428

$ws.Range("C9").Value = -0.4038821969541538
$ws.Range("D9").Value = 0.6901984459789681

$ws.Range("C10").Value = -0.6409403460663854
$ws.Range("D10").Value = 0.5281822797995188

$ws.Range("C11").Value = -0.4510301530130711
$ws.Range("D11").Value = 0.6563820532619253
